# Added SSO authentication for FCV tests:
# update the "build" value (column A) for the rows whose build was still
# marked as the older "test" placeholder so they report the newer build
# ("33c392bc2b built at 2020-09-17 13:46") already used by the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newBuild = "33c392bc2b built at 2020-09-17 13:46`n"

$rows = @(5, 10, 13, 19, 22)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = $newBuild
    $ws.Rows.Item($r).AutoFit()
}
